$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.200.06"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.587.37"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.05"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("E10").Value = "  -1.76%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +0.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.811.83"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.590.43"
$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("E14").Value = "  -1.56%  "

$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.94"
$ws.Range("E16").Value = "  -1.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.206.56"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0724"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "214.41"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("E20").Value = "  -2.62%  "

$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("E22").Value = "  -0.86%  "

$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("E25").Value = "  -1.17%  "

$ws.Range("E26").Value = "  -0.10%  "

$ws.Range("E27").Value = "  -0.79%  "

$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.12"
$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("E30").Value = "  -1.75%  "

$ws.Range("E31").Value = "  +0.62%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.390.00"
$ws.Range("E33").Value = "  +6.78%  "

$ws.Range("E34").Value = "  -1.67%  "

$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("E36").Value = "  -1.33%  "

$ws.Range("E37").Value = "  -5.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.817"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.940"
$ws.Range("E42").Value = "  -15.11%  "

$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.723.72"
$ws.Range("E45").Value = "  -0.23%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.97"
$ws.Range("E46").Value = "  -2.80%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.10"
$ws.Range("E47").Value = "  -2.76%  "

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.27%  "
